# Apply the updates described by the diff:
# 1. Column C (rows 2-10): update the "Förändrad" date serial from 45207 to 45208.
# 2. Rows 2-4, columns S, T, V, W, X, Y: update the HYPERLINK formulas so the
#    path segment "Logging_OSTERSUND" becomes "Logging_2380".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update column C values for rows 2 through 10 ---
for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)   # Column C
    if ($cell.Value2 -eq 45207) {
        $cell.Value = 45208
    }
}

# --- 2. Update hyperlink formulas in rows 2-4 for columns S, T, V, W, X, Y ---
$linkColumns = @("S", "T", "V", "W", "X", "Y")
for ($row = 2; $row -le 4; $row++) {
    foreach ($col in $linkColumns) {
        $addr = "$col$row"
        $cell = $ws.Range($addr)
        $formula = $cell.Formula
        if ($formula -and $formula.Length -gt 0) {
            $newFormula = $formula.Replace("Logging_OSTERSUND", "Logging_2380")
            if ($newFormula -ne $formula) {
                $cell.Formula = $newFormula
            }
        }
    }
}
